# Updates the Price (D) and Volume(1h) (E) columns of the cryptos
# list sheet to the latest scraped values, preserving each cell
# as text (matching the workbook's existing inline-string cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.856.03"
$ws.Range("E2").Value = "  +0.52%  "

$ws.Range("D3").Value = "3.105.83"
$ws.Range("E3").Value = "  +3.94%  "

$ws.Range("E4").Value = "  +0.01%  "

$d = $ws.Range("D5")
$d.NumberFormat = "@"
$d.Value = "390.36"
$d.ClearFormats()
$ws.Range("E5").Value = "  +2.30%  "

$d = $ws.Range("D6")
$d.NumberFormat = "@"
$d.Value = "103.84"
$d.ClearFormats()
$ws.Range("E6").Value = "  -0.44%  "

$d = $ws.Range("D7")
$d.NumberFormat = "@"
$d.Value = "0.546"
$d.ClearFormats()
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  +0.01%  "

$d = $ws.Range("D9")
$d.NumberFormat = "@"
$d.Value = "0.593"
$d.ClearFormats()
$ws.Range("E9").Value = "  -0.44%  "

$d = $ws.Range("D10")
$d.NumberFormat = "@"
$d.Value = "37.24"
$d.ClearFormats()
$ws.Range("E10").Value = "  +1.32%  "

$ws.Range("E11").Value = "  -0.13%  "

$d = $ws.Range("D12")
$d.NumberFormat = "@"
$d.Value = "0.0864"
$d.ClearFormats()
$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("D13").Value = "3.599.63"
$ws.Range("E13").Value = "  +3.90%  "

$d = $ws.Range("D14")
$d.NumberFormat = "@"
$d.Value = "18.77"
$d.ClearFormats()
$ws.Range("E14").Value = "  +1.37%  "

$ws.Range("E15").Value = "  +0.59%  "

$ws.Range("D16").Value = "3.095.87"
$ws.Range("E16").Value = "  +3.68%  "

$d = $ws.Range("D17")
$d.NumberFormat = "@"
$d.Value = "0.987"
$d.ClearFormats()
$ws.Range("E17").Value = "  -0.85%  "

$d = $ws.Range("D18")
$d.NumberFormat = "@"
$d.Value = "10.80"
$d.ClearFormats()
$ws.Range("E18").Value = "  -3.62%  "

$ws.Range("D19").Value = "51.919.53"
$ws.Range("E19").Value = "  +0.59%  "

$ws.Range("E20").Value = "  +3.94%  "

$ws.Range("E21").Value = "  -0.70%  "

$ws.Range("D22").Value = "0.0₃0969"
$ws.Range("E22").Value = "  +0.58%  "

$d = $ws.Range("D23")
$d.NumberFormat = "@"
$d.Value = "70.10"
$d.ClearFormats()
$ws.Range("E23").Value = "  -0.57%  "

$d = $ws.Range("D24")
$d.NumberFormat = "@"
$d.Value = "268.32"
$d.ClearFormats()
$ws.Range("E24").Value = "  +0.34%  "

$d = $ws.Range("D25")
$d.NumberFormat = "@"
$d.Value = "3.14"
$d.ClearFormats()
$ws.Range("E25").Value = "  -2.70%  "

$d = $ws.Range("D26")
$d.NumberFormat = "@"
$d.Value = "8.15"
$d.ClearFormats()
$ws.Range("E26").Value = "  +3.20%  "

$d = $ws.Range("D27")
$d.NumberFormat = "@"
$d.Value = "27.22"
$d.ClearFormats()
$ws.Range("E27").Value = "  +4.15%  "

$ws.Range("E28").Value = "  +0.15%  "

$d = $ws.Range("D29")
$d.NumberFormat = "@"
$d.Value = "7.21"
$d.ClearFormats()
$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("E30").Value = "  +0.13%  "

$ws.Range("E31").Value = "  -0.82%  "

$ws.Range("E32").Value = "  -0.18%  "

$d = $ws.Range("D33")
$d.NumberFormat = "@"
$d.Value = "35.55"
$d.ClearFormats()
$ws.Range("E33").Value = "  +2.89%  "

$ws.Range("E34").Value = "  +0.69%  "

$d = $ws.Range("D35")
$d.NumberFormat = "@"
$d.Value = "50.41"
$d.ClearFormats()
$ws.Range("E35").Value = "  -1.78%  "

$d = $ws.Range("D36")
$d.NumberFormat = "@"
$d.Value = "0.0451"
$d.ClearFormats()
$ws.Range("E36").Value = "  +1.11%  "

$ws.Range("E37").Value = "  -0.15%  "

$ws.Range("E38").Value = "  +3.67%  "

$ws.Range("E39").Value = "  +7.11%  "

$ws.Range("E40").Value = "  +2.67%  "

$ws.Range("E41").Value = "  +0.60%  "

$d = $ws.Range("D42")
$d.NumberFormat = "@"
$d.Value = "16.92"
$d.ClearFormats()
$ws.Range("E42").Value = "  -0.20%  "

$d = $ws.Range("D43")
$d.NumberFormat = "@"
$d.Value = "129.05"
$d.ClearFormats()
$ws.Range("E43").Value = "  +1.22%  "

$ws.Range("E45").Value = "  -3.68%  "

$d = $ws.Range("D46")
$d.NumberFormat = "@"
$d.Value = "22.31"
$d.ClearFormats()
$ws.Range("E46").Value = "  +3.83%  "

$d = $ws.Range("D47")
$d.NumberFormat = "@"
$d.Value = "2.50"
$d.ClearFormats()
$ws.Range("E47").Value = "  +6.38%  "

$ws.Range("E48").Value = "  +2.15%  "

$ws.Range("D49").Value = "2.052.91"
$ws.Range("E49").Value = "  +0.88%  "

$ws.Range("D50").Value = "3.415.50"

$ws.Range("E51").Value = "  -1.00%  "
